# Enabling all testcases KFP
# Set Runmode (column E) to "Yes" for all testcase rows except
# row 19 (TC19_Verify_PunchOut_User) and row 26 (TC30_Verify_pagination_SortBy_filteronPLP),
# which remain "No".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

$rowsToEnable = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,20,21,22,23,24,25,27,28)

foreach ($r in $rowsToEnable) {
    $ws.Range("E$r").Value = "Yes"
}

# Update the selected cell to match the saved view state.
$ws.Range("E19").Select()
